$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E24").Value = 226.31
$ws.Range("E26").Value = 1445.53
